{"js": "// Replace each two-digit-by-two-digit multiplication problem in the\n// document with its new equation, per the commit's regenerated answer key.\n// Every occurrence of the old equation text is unique in the document, so a\n// simple exact-text search/replace (matchCase + wholeWords off, since the\n// strings contain \"\u00d7\" and \"=\") is sufficient and robust regardless of how\n// the table rows/cells are indexed.\n\nconst replacements = [\n  [\"35\u00d740=1400\", \"34\u00d713=442\"],\n  [\"43\u00d799=4257\", \"57\u00d715=855\"],\n  [\"49\u00d795=4655\", \"51\u00d730=1530\"],\n  [\"82\u00d774=6068\", \"17\u00d762=1054\"],\n  [\"94\u00d726=2444\", \"29\u00d718=522\"],\n  [\"54\u00d735=1890\", \"36\u00d732=1152\"],\n  [\"43\u00d776=3268\", \"21\u00d769=1449\"],\n  [\"87\u00d757=4959\", \"55\u00d770=3850\"],\n  [\"13\u00d791=1183\", \"51\u00d742=2142\"],\n  [\"92\u00d790=8280\", \"27\u00d725=675\"],\n  [\"99\u00d789=8811\", \"49\u00d784=4116\"],\n  [\"80\u00d742=3360\", \"18\u00d795=1710\"],\n  [\"38\u00d741=1558\", \"21\u00d779=1659\"],\n  [\"15\u00d726=390\", \"74\u00d792=6808\"],\n  [\"26\u00d786=2236\", \"38\u00d761=2318\"],\n  [\"44\u00d726=1144\", \"60\u00d753=3180\"],\n  [\"55\u00d718=990\", \"77\u00d723=1771\"],\n  [\"94\u00d749=4606\", \"24\u00d747=1128\"],\n  [\"61\u00d777=4697\", \"24\u00d737=888\"],\n  [\"60\u00d794=5640\", \"56\u00d782=4592\"],\n  [\"67\u00d794=6298\", \"37\u00d711=407\"],\n  [\"80\u00d741=3280\", \"81\u00d777=6237\"],\n  [\"68\u00d773=4964\", \"24\u00d754=1296\"],\n  [\"31\u00d777=2387\", \"66\u00d766=4356\"],\n  [\"96\u00d752=4992\", \"77\u00d793=7161\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Could not find text to replace: \"${oldText}\"`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each two-digit-by-two-digit multiplication problem in the\n# document with its new equation, per the commit's regenerated answer key.\n# Every occurrence of the old equation text is unique in the document, so a\n# simple Find/Replace-All per pair is sufficient and robust regardless of\n# how the table rows/cells are indexed.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"35\u00d740=1400\", \"34\u00d713=442\"),\n    @(\"43\u00d799=4257\", \"57\u00d715=855\"),\n    @(\"49\u00d795=4655\", \"51\u00d730=1530\"),\n    @(\"82\u00d774=6068\", \"17\u00d762=1054\"),\n    @(\"94\u00d726=2444\", \"29\u00d718=522\"),\n    @(\"54\u00d735=1890\", \"36\u00d732=1152\"),\n    @(\"43\u00d776=3268\", \"21\u00d769=1449\"),\n    @(\"87\u00d757=4959\", \"55\u00d770=3850\"),\n    @(\"13\u00d791=1183\", \"51\u00d742=2142\"),\n    @(\"92\u00d790=8280\", \"27\u00d725=675\"),\n    @(\"99\u00d789=8811\", \"49\u00d784=4116\"),\n    @(\"80\u00d742=3360\", \"18\u00d795=1710\"),\n    @(\"38\u00d741=1558\", \"21\u00d779=1659\"),\n    @(\"15\u00d726=390\", \"74\u00d792=6808\"),\n    @(\"26\u00d786=2236\", \"38\u00d761=2318\"),\n    @(\"44\u00d726=1144\", \"60\u00d753=3180\"),\n    @(\"55\u00d718=990\", \"77\u00d723=1771\"),\n    @(\"94\u00d749=4606\", \"24\u00d747=1128\"),\n    @(\"61\u00d777=4697\", \"24\u00d737=888\"),\n    @(\"60\u00d794=5640\", \"56\u00d782=4592\"),\n    @(\"67\u00d794=6298\", \"37\u00d711=407\"),\n    @(\"80\u00d741=3280\", \"81\u00d777=6237\"),\n    @(\"68\u00d773=4964\", \"24\u00d754=1296\"),\n    @(\"31\u00d777=2387\", \"66\u00d766=4356\"),\n    @(\"96\u00d752=4992\", \"77\u00d793=7161\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $found = $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n\n    if (-not $found) {\n        throw \"Could not find text to replace: $oldText\"\n    }\n}\n"}
